# "add tabel format baru" - renumber the table headings (new tables were
# inserted earlier in the series, bumping Tabel 4.2.3/4.2.4/4.2.5 to
# 4.2.5/4.2.6/4.2.7) and roll the reporting year from 2020 to 2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1 (columns A-M): "Tabel 4.2.3" -> "Tabel 4.2.5" -------------
$ws.Range("H1").Value = "Tabel 4.2.5"
$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan Wundulako. 2021"
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan Wundulako, 2021"
$ws.Range("B2").Value = "Number of Medical Personnel by Kelurahan/ Village in Wundulako Subdistrict, 2021"
$ws.Range("I2").Value = "Number of Immunized Babies by Types of Immunization and Kelurahan/Village Wundulako Subdistrict, 2021"

# --- Block 2 (columns P-S): "Tabel 4.2.4." -> "Tabel 4.2.6." -----------
# Use Characters() on just the numeric run so the "Tabel" / " 4.2.6."
# rich-text split is preserved instead of collapsing the whole cell.
$ws.Range("P1").Characters(7, 6).Text = "4.2.6."
$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan Wundulako, 2021"
$ws.Range("Q2").Value = "Number of Woman Giving Brth and Birth Assisted by Paramedics by Kelurahan/Village in Wundulako Subdistrict, 2021"

# --- Block 3 (columns W-Z): "Tabel 4.2.5." -> "Tabel 4.2.7." -----------
$ws.Range("W1").Characters(7, 6).Text = "4.2.7."
$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan Wundulako, 2021"
$ws.Range("X2").Value = "Number of Fertile Age Couples and Family Planning Members by Kelurahan/Village in Wundulako Subdistrict, 2021"
